# Corrected excel sheets for application fix issues
$wb = $excel.ActiveWorkbook

# --- Summary sheet: Fees row (row 4) updated from 50 to 100, and selection moved ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A4").Value = 100
$wsSummary.Range("B4").Value = 100

# --- Repayment Schedule sheet: disbursement-row fee figures updated from 50 to 100 ---
$wsRepay = $wb.Worksheets.Item("Repayment Schedule")
$wsRepay.Range("I2").Value = 100
$wsRepay.Range("K2").Value = 100
$wsRepay.Range("L2").Value = 100

# --- Transactions sheet: fee figures updated from 50 to 100 ---
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("E2").Value = 100
$wsTrans.Range("H2").Value = 100

# --- Update each sheet's cursor/selection to match the saved view state ---
$wsSummary.Select()
$wsSummary.Range("G15:G16").Select()

$wsRepay.Select()
$wsRepay.Range("K19:L19").Select()

# Transactions is the sheet left active/selected when the workbook was saved
$wsTrans.Select()
$wsTrans.Range("K9").Select()
